# 4.2.2.1a.xlsx — add the 2022 data column (M) to the coverage table.
#
# Column L (2021) is the rightmost existing year; we append a new column M
# for 2022 using the same per-row formatting as the existing data, then
# write in the 2022 figures. Row 14 ("Ysyk-Kul oblast" total row) uses a
# style that doesn't exist yet elsewhere in the sheet, so it is built from
# the nearby bold-row style plus the table's numeric (0.0) format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# Each entry: new M-column row, a same-styled source cell to copy formatting
# from, and the 2022 value to write.
$data = @(
    @{ Row = 4;  Src = "K4";  Val = 2022 }
    @{ Row = 5;  Src = "K5";  Val = 24.6 }
    @{ Row = 6;  Src = "K6";  Val = 40.7 }
    @{ Row = 7;  Src = "K6";  Val = 20.7 }
    @{ Row = 8;  Src = "K5";  Val = 26.6 }
    @{ Row = 9;  Src = "K6";  Val = 44.5 }
    @{ Row = 10; Src = "K6";  Val = 21.9 }
    @{ Row = 11; Src = "K5";  Val = 21.9 }
    @{ Row = 12; Src = "K6";  Val = 35.3 }
    @{ Row = 13; Src = "K6";  Val = 17.6 }
    @{ Row = 14; Src = "K5";  Val = 28;   NumFmt = "0.0" }
    @{ Row = 15; Src = "K6";  Val = 44.9 }
    @{ Row = 16; Src = "K6";  Val = 21.5 }
    @{ Row = 17; Src = "K5";  Val = 36.2 }
    @{ Row = 18; Src = "K6";  Val = 53.1 }
    @{ Row = 19; Src = "K6";  Val = 33.4 }
    @{ Row = 20; Src = "K5";  Val = 20.2 }
    @{ Row = 21; Src = "K6";  Val = 15.4 }
    @{ Row = 22; Src = "K6";  Val = 20.5 }
    @{ Row = 23; Src = "K5";  Val = 27.1 }
    @{ Row = 24; Src = "K6";  Val = 36.1 }
    @{ Row = 25; Src = "K6";  Val = 25.2 }
    @{ Row = 26; Src = "K5";  Val = 24.2 }
    @{ Row = 27; Src = "K6";  Val = 46.5 }
    @{ Row = 28; Src = "K6";  Val = 20.3 }
    @{ Row = 29; Src = "K5";  Val = 40.5 }
    @{ Row = 30; Src = "L30"; Val = 44.5 }
)

foreach ($item in $data) {
    $target = $ws.Range("M" + $item.Row)

    # Copy number/font/border/fill formatting from an existing cell that
    # already carries the style this row needs.
    $ws.Range($item.Src).Copy()
    $target.PasteSpecial($xlPasteFormats)

    if ($item.ContainsKey("NumFmt")) {
        $target.NumberFormat = $item.NumFmt
    }

    $target.Value = $item.Val
}

$excel.CutCopyMode = $false

# Match the author's final selection.
$ws.Range("N7").Select() | Out-Null
